# Add a new "Save" column (H) to the s_vals sheet, matching the
# header formatting used by the existing columns and a default 0
# value in the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in H1, using the same style as the other header
# cells (e.g. G1 "sum") so it matches the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell in H2 for the "Save" column.
$ws.Range("H2").Value = 0
